$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the unified "DataNode" concept
# (was "Property1").
$ws.Name = "DataNode"

# Move/restore the active selection onto D26 (matches the saved view state).
$ws.Range("D26").Select()
